# Append the new run-log row (row 60) to the Nalco run log sheet, matching
# the formatting of the previous row (row 59) and the data from the latest
# run (2025-08-26 09:40:54 UTC / SKIPPED - no PDF change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the style of the last existing data row (row 59) onto the new row
# 60 so every cell (including the normally-blank "Saved PDF" / "Total Rows
# After" columns) carries the same centered style used throughout the log.
$ws.Range("A59:H59").Copy()
$ws.Range("A60:H60").PasteSpecial(-4122)

# Fill in the new run's data.
$ws.Range("A60").Value = "2025-08-26 09:40:54 UTC"
$ws.Range("B60").Value = "2025-08-26 15:10:54 IST"
$ws.Range("C60").Value = "SKIPPED"
$ws.Range("D60").Value = "No change in PDF. Skipping download & Excel update."
$ws.Range("E60").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Range("G60").Value = 0
